$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E8").Value = "aaron_clark.png"
$ws.Range("D8").Value = "Workstream Co-lead"
